{"js": "// The document contains four occurrences of an \"<id>p043r_aN</id>\" marker,\n// each split across three runs: \"<id>\" (Courier New / #7f6000 / 9pt),\n// \"p043r_aN\" (plain, black), \"</id>\" (Courier New / #7f6000 / 9pt).\n// The edit collapses each trio into a single run reading\n// \"<id>p043r_N</id>\" (dropping the \"a\" before the digit) that keeps the\n// formatting of the original \"<id>\" run.\nconst mapping = [\n  [\"<id>p043r_a1</id>\", \"<id>p043r_1</id>\"],\n  [\"<id>p043r_a2</id>\", \"<id>p043r_2</id>\"],\n  [\"<id>p043r_a3</id>\", \"<id>p043r_3</id>\"],\n  [\"<id>p043r_a4</id>\", \"<id>p043r_4</id>\"],\n];\n\nfor (const [oldText, newText] of mapping) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    // insertText with Replace collapses the whole matched range (which may\n    // span multiple runs) into a single run, keeping the formatting of the\n    // range's (first) run \u2014 i.e. the Courier New / #7f6000 \"<id>\" styling.\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains four occurrences of an \"<id>p043r_aN</id>\" marker,\n# each split across three runs: \"<id>\" (Courier New / #7f6000 / 9pt),\n# \"p043r_aN\" (plain, black), \"</id>\" (Courier New / #7f6000 / 9pt).\n# The edit collapses each trio into a single run reading\n# \"<id>p043r_N</id>\" (dropping the \"a\" before the digit) that keeps the\n# formatting of the original \"<id>\" run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = \"<id>p043r_a1</id>\"; New = \"<id>p043r_1</id>\" },\n    @{ Old = \"<id>p043r_a2</id>\"; New = \"<id>p043r_2</id>\" },\n    @{ Old = \"<id>p043r_a3</id>\"; New = \"<id>p043r_3</id>\" },\n    @{ Old = \"<id>p043r_a4</id>\"; New = \"<id>p043r_4</id>\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Find the whole \"<id>...</id>\" span (it reads across the three runs)\n    # and replace it in one shot; Word collapses the matched range into a\n    # single run carrying the formatting of the range's leading run.\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
